$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" list had two rows (24: LOM3254 lab indication, 25: LOB1053
# physics requirement). The build re-ordered the two entries so that the
# LOB1053 requirement now comes first (row 24) and the LOM3254 indication
# comes second (row 25). Swap the values of the two rows (columns B and C).

$reqB24 = $ws.Range("B24").Value()
$reqC24 = $ws.Range("C24").Value()
$reqB25 = $ws.Range("B25").Value()
$reqC25 = $ws.Range("C25").Value()

$ws.Range("B24").Value = $reqB25
$ws.Range("C24").Value = $reqC25
$ws.Range("B25").Value = $reqB24
$ws.Range("C25").Value = $reqC24
